# Update: Threat Alert Report - 2026-01-30 06:29
#
# The AHB threats table (rows 2-18) is replaced by a refreshed 13-row
# report (rows 2-14): dates roll forward, fares/diffs are recomputed, the
# matched competitor flight/airline shifts, and the MEDIUM THREAT flag
# moves from row 10 (old) to row 9 (new). The trailing four rows (old
# rows 15-18) are removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the four trailing data rows (old rows 15-18) -----------------
# Doing this first means rows 2-14 keep their row numbers while we fill in
# the refreshed data below, and the sheet's used-range/dimension shrinks to
# A1:K14 automatically.
$ws.Rows("15:18").Delete()

# --- 2. Refreshed data for rows 2-14 ---------------------------------------
$rows = @(
    @{ Row=2;  Date="14-FEB-26"; Flight="SM-443"; Airline="Nile Air NP-243";         OALFare=7651;  OurFare=7914;  FareDif=-263;  Impact="LOW THREAT" },
    @{ Row=3;  Date="19-FEB-26"; Flight="SM-443"; Airline="Nile Air NP-143";         OALFare=7651;  OurFare=7914;  FareDif=-263;  Impact="LOW THREAT" },
    @{ Row=4;  Date="21-FEB-26"; Flight="SM-443"; Airline="Nile Air NP-243";         OALFare=7651;  OurFare=7914;  FareDif=-263;  Impact="LOW THREAT" },
    @{ Row=5;  Date="28-FEB-26"; Flight="SM-443"; Airline="Nile Air NP-243";         OALFare=7651;  OurFare=7914;  FareDif=-263;  Impact="LOW THREAT" },
    @{ Row=6;  Date="21-MAR-26"; Flight="SM-443"; Airline="Nile Air NP-243";         OALFare=7651;  OurFare=7914;  FareDif=-263;  Impact="LOW THREAT" },
    @{ Row=7;  Date="26-MAR-26"; Flight="SM-443"; Airline="Nile Air NP-143";         OALFare=13672; OurFare=14575; FareDif=-903;  Impact="LOW THREAT" },
    @{ Row=8;  Date="28-MAR-26"; Flight="SM-443"; Airline="Nile Air NP-243";         OALFare=18702; OurFare=20220; FareDif=-1518; Impact="LOW THREAT" },
    @{ Row=9;  Date="30-MAR-26"; Flight="SM-443"; Airline="Air Arabia Egypt E5-511"; OALFare=7170;  OurFare=10160; FareDif=-2990; Impact="MEDIUM THREAT - MONITOR" },
    @{ Row=10; Date="01-APR-26"; Flight="SM-443"; Airline="Nile Air NP-143";         OALFare=8379;  OurFare=8818;  FareDif=-439;  Impact="LOW THREAT" },
    @{ Row=11; Date="04-APR-26"; Flight="SM-443"; Airline="Air Arabia Egypt E5-513"; OALFare=7240;  OurFare=7488;  FareDif=-248;  Impact="LOW THREAT" },
    @{ Row=12; Date="11-MAY-26"; Flight="SM-443"; Airline="Air Arabia Egypt E5-511"; OALFare=7170;  OurFare=8328;  FareDif=-1158; Impact="LOW THREAT" },
    @{ Row=13; Date="16-MAY-26"; Flight="SM-443"; Airline="Air Arabia Egypt E5-513"; OALFare=7240;  OurFare=8328;  FareDif=-1088; Impact="LOW THREAT" },
    @{ Row=14; Date="18-MAY-26"; Flight="SM-443"; Airline="Air Arabia Egypt E5-511"; OALFare=7170;  OurFare=8328;  FareDif=-1158; Impact="LOW THREAT" }
)

# Scratch cell used to stage each date string as genuine text (formatted
# "@") before copying just the VALUE onto column A - this stops Excel's
# auto date-recognition from turning "14-FEB-26" into a date serial while
# leaving the destination cell's existing style (border/alignment) intact.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

foreach ($r in $rows) {
    $scratch.Value = $r.Date
    $scratch.Copy()
    $ws.Cells.Item($r.Row, 1).PasteSpecial("xlPasteValues")

    $ws.Cells.Item($r.Row, 2).Value = $r.Flight
    $ws.Cells.Item($r.Row, 3).Value = $r.Airline
    $ws.Cells.Item($r.Row, 4).Value = $r.OALFare
    $ws.Cells.Item($r.Row, 5).Value = $r.OurFare
    $ws.Cells.Item($r.Row, 6).Value = $r.FareDif
    $ws.Cells.Item($r.Row, 10).Value = $r.Impact
}

$scratch.Clear()

# --- 3. IMPACT cell styling: the MEDIUM THREAT highlight moves from the
#        old row 10 to the new row 9, and row 10 reverts to the regular
#        LOW THREAT look. The text we just wrote into J9/J10 already
#        matches what used to sit in J10/J11 respectively, so copying
#        those whole cells (value + style together) over reproduces both
#        the correct text AND the correct highlight style index in one
#        step - Range.Copy(destination) carries the source's style, unlike
#        PasteSpecial in this host which leaves the destination's existing
#        style untouched.
$ws.Range("J10").Copy($ws.Range("J9"))
$ws.Range("J11").Copy($ws.Range("J10"))

Write-Output "Threat table refreshed: rows 2-14 updated, rows 15-18 removed."
